# Update forecast summary workbook:
#  - Insert a new "Week_Start_Date" column (B) into the "Forecast Comparison" sheet
#  - Normalize the "Week" labels (W01 -> W1, etc.)
#  - Refresh the "MyForecast" numbers with the corrected output
#  - Store "is_holiday_week" as a boolean instead of a 0/1 number
#  - Refresh the derived totals on the "Summary" sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")
$summary = $wb.Worksheets.Item("Summary")

# --- Insert the new "Week_Start_Date" column as column B -------------------
# (this shifts ASIN/MyForecast/Amazon*/Product Title/is_holiday_week one
# column to the right: B->C, C->D, D->E, E->F, F->G, G->H, H->I, I->J)
$ws.Columns.Item(2).Insert()
$ws.Range("B1").Value = "Week_Start_Date"

# --- Fill in the Week_Start_Date values (kept as literal text, not dates) --
$weekStartDates = @{
    2  = "2025-01-05"
    3  = "2025-01-12"
    4  = "2025-01-19"
    5  = "2025-01-26"
    6  = "2025-02-02"
    7  = "2025-02-09"
    8  = "2025-02-16"
    9  = "2025-02-23"
    10 = "2025-03-02"
    11 = "2025-03-09"
    12 = "2025-03-16"
    13 = "2025-03-23"
    14 = "2025-03-30"
    15 = "2025-04-06"
    16 = "2025-04-13"
    17 = "2025-04-20"
}

# Force text formatting on the whole fill range up front so Excel does not
# auto-convert the "YYYY-MM-DD" strings into date serials, then clear the
# formatting back off again afterwards so no cell-level style is left behind.
$dateRange = $ws.Range("B2:B17")
$dateRange.NumberFormat = "@"
foreach ($row in 2..17) {
    $ws.Range("B$row").Value = $weekStartDates[$row]
}
$dateRange.ClearFormats()

# --- Normalize the Week labels (strip the leading zero) --------------------
$weekLabels = @{
    2  = "W1"
    3  = "W2"
    4  = "W3"
    5  = "W4"
    6  = "W5"
    7  = "W6"
    8  = "W7"
    9  = "W8"
    10 = "W9"
}
foreach ($row in $weekLabels.Keys) {
    $ws.Range("A$row").Value = $weekLabels[$row]
}

# --- Refresh the MyForecast column (now column D) with corrected output ----
$myForecast = @{
    2  = 171
    3  = 182
    4  = 188
    5  = 203
    6  = 167
    7  = 154
    8  = 146
    9  = 148
    10 = 156
    11 = 163
    12 = 185
    13 = 180
    14 = 170
    15 = 158
    16 = 147
    17 = 151
}
foreach ($row in 2..17) {
    $ws.Range("D$row").Value = $myForecast[$row]
}

# --- Store is_holiday_week (now column J) as a real boolean ----------------
foreach ($row in 2..17) {
    $ws.Range("J$row").Value = $false
}

# --- Refresh the derived totals on the Summary sheet ------------------------
# (these are stored as text in the workbook, not numbers, so force text
# formatting the same way as the week-start dates above)
$summaryRange = $summary.Range("B9:B12")
$summaryRange.NumberFormat = "@"
$summary.Range("B9").Value = "2668"
$summary.Range("B10").Value = "1358"
$summary.Range("B11").Value = "743"
$summary.Range("B12").Value = "203"
$summaryRange.ClearFormats()
